$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# Add the "Note:" label and its text in row 11, columns F and G.
$ws.Range("F11").Value = "Note:"
$ws.Range("F11").Style = "Accent1"

$ws.Range("G11").Value = "Not a Microsoft supported deployment topology"
$ws.Range("G11").Interior.Color = 65535

# Update the selection to match the new active cell/range.
$ws.Range("F11:G11").Select()
